## Deploying to gh-pages from @ LinuxForHealth/alvearie-fhir-ig
## Update the StructureDefinition metadata (URL, Version, Date, Publisher)
## to reflect the move from IBM/Alvearie to the LinuxForHealth project,
## and drop the stray duplicate FHIR constraint text that used to be
## copied onto the top-level "Extension" row of the Elements sheet.

$wb = $excel.ActiveWorkbook

# --- Metadata sheet -------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/episode-summary-group-code"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

# --- Elements sheet ---------------------------------------------------
$elements = $wb.Worksheets.Item("Elements")

# The Extension.url row fixes its value to this StructureDefinition's own
# canonical URL, so it must track the same URL update as the Metadata
# sheet.
$elements.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/episode-summary-group-code"

# The "Extension" (root) row had the ele-1/ext-1 constraint text
# erroneously duplicated into its Constraint(s) column (AI2). That
# constraint really only belongs on the child rows (e.g. AI4, AI6), so
# clear it from the parent row.
$elements.Range("AI2").Value = ""
